# Update latest output (run 52)

$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule" ---
$schedule = $wb.Worksheets.Item("Schedule")
$schedule.Range("E2").Value = 714.2583240000001
$schedule.Range("F2").Value = 11.80982678571429

# --- Sheet "Detailed" ---
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Range("B9").Value = 36.06
$detailed.Range("B10").Value = 36.06

$detailed.Range("B11").Value = 57.3
$detailed.Range("C11").Value = "historical"

$detailed.Range("C12").Value = "historical"

$detailed.Range("B17").Value = 0.7
$detailed.Range("B18").Value = 0.7
$detailed.Range("B19").Value = 33.20785
$detailed.Range("B20").Value = 36.06054

$detailed.Range("B22").Value = 36.06011

$detailed.Range("B24").Value = 43.62934

$detailed.Range("B30").Value = 31.00384
$detailed.Range("B31").Value = 25.73729
$detailed.Range("B32").Value = 29.57574
$detailed.Range("B33").Value = 33.21128
$detailed.Range("B34").Value = 33.26438
$detailed.Range("B35").Value = 8.411429999999999
$detailed.Range("B36").Value = -3.09877
$detailed.Range("B37").Value = -3.01756
$detailed.Range("B38").Value = -2.91645
$detailed.Range("B39").Value = -2.86201
$detailed.Range("B40").Value = 3.42641
$detailed.Range("B41").Value = 9.65208
$detailed.Range("B42").Value = 9.71571
$detailed.Range("B43").Value = 19.54547
$detailed.Range("B44").Value = 8.333170000000001
$detailed.Range("B45").Value = 6.5731
$detailed.Range("B46").Value = 30.19732
$detailed.Range("B47").Value = 57.06008

$detailed.Range("B49").Value = 56.97995
